$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 4818.3335
$ws.Range("I12").Value = 2700
$ws.Range("J12").Value = 6936.6665
$ws.Range("K12").Value = 2700
$ws.Range("L12").Value = 6936.6665
$ws.Range("M12").Value = -2530
$ws.Range("N12").Value = -7276.6665
$ws.Range("H32").Value = 2454.6667
$ws.Range("I32").Value = 2649.8
$ws.Range("J32").Value = 2357.1
$ws.Range("K32").Value = 2649.8
$ws.Range("L32").Value = 2357.1
$ws.Range("M32").Value = -2323.8
$ws.Range("N32").Value = -3009.1
$ws.Range("H58").Value = 4897.3335
$ws.Range("I58").Value = 573.125
$ws.Range("J58").Value = 9839.286
$ws.Range("K58").Value = 1719.375
$ws.Range("L58").Value = 29517.858
$ws.Range("M58").Value = -1569.375
$ws.Range("N58").Value = -29817.858
$ws.Range("H87").Value = 74518.17999999999
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 74518.17999999999
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 74518.17999999999
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -77014.17999999999
$ws.Range("H90").Value = 74518.17999999999
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 74518.17999999999
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 223554.54
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -236034.54
$ws.Range("H136").Value = 50000
$ws.Range("J136").Value = 50000
$ws.Range("L136").Value = 50000
$ws.Range("N136").Value = -60200
$ws.Range("H137").Value = 4423.143
$ws.Range("I137").Value = 4419.7856
$ws.Range("K137").Value = 13259.3568
$ws.Range("M137").Value = -10709.3568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3352
$ws.Range("I32").Value = 3259.2075
$ws.Range("K32").Value = 3259.2075
$ws.Range("M32").Value = -2972.2075
$ws.Range("H61").Value = 3559
$ws.Range("I61").Value = 1337
$ws.Range("K61").Value = 1337
$ws.Range("M61").Value = -1125
$ws.Range("H136").Value = 3559
$ws.Range("I136").Value = 1337
$ws.Range("K136").Value = 4011
$ws.Range("M136").Value = -1461

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 70000
$ws.Range("J52").Value = 70000
$ws.Range("L52").Value = 70000
$ws.Range("N52").Value = -70526
$ws.Range("H86").Value = 2930.5715
$ws.Range("I86").Value = 1266.6666
$ws.Range("J86").Value = 4178.5
$ws.Range("K86").Value = 1266.6666
$ws.Range("L86").Value = 4178.5
$ws.Range("M86").Value = -143.6666
$ws.Range("N86").Value = -6424.5
$ws.Range("H89").Value = 2930.5715
$ws.Range("I89").Value = 1266.6666
$ws.Range("J89").Value = 4178.5
$ws.Range("K89").Value = 6333.333000000001
$ws.Range("L89").Value = 20892.5
$ws.Range("M89").Value = -717.3330000000005
$ws.Range("N89").Value = -32124.5
$ws.Range("H105").Value = 47677.773
$ws.Range("J105").Value = 2302.2666
$ws.Range("L105").Value = 2302.2666
$ws.Range("N105").Value = -5796.2666
$ws.Range("H121").Value = 70000
$ws.Range("J121").Value = 70000
$ws.Range("L121").Value = 70000
$ws.Range("N121").Value = -73494
$ws.Range("H128").Value = 7833.3335
$ws.Range("I128").Value = 7833.3335
$ws.Range("K128").Value = 23500.0005
$ws.Range("M128").Value = -21010.0005
$ws.Range("H134").Value = 2828.606
$ws.Range("I134").Value = 1941.9642
$ws.Range("K134").Value = 5825.892599999999
$ws.Range("M134").Value = -3290.892599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3668.8484
$ws.Range("I132").Value = 2836.2
$ws.Range("K132").Value = 8508.599999999999
$ws.Range("M132").Value = -5978.599999999999
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4276249.5
$ws.Range("I4").Value = 2556110
$ws.Range("K4").Value = 7668330
$ws.Range("M4").Value = -7668218
$ws.Range("H33").Value = 3086431
$ws.Range("I33").Value = 3086431
$ws.Range("K33").Value = 18518586
$ws.Range("M33").Value = -18518303
$ws.Range("H40").Value = 237.89473
$ws.Range("I40").Value = 136.42857
$ws.Range("J40").Value = 522
$ws.Range("K40").Value = 545.71428
$ws.Range("L40").Value = 2088
$ws.Range("M40").Value = -476.71428
$ws.Range("N40").Value = -2226
$ws.Range("H44").Value = 2434.3333
$ws.Range("I44").Value = 303
$ws.Range("K44").Value = 909
$ws.Range("M44").Value = -511
$ws.Range("H68").Value = 1883.1666
$ws.Range("I68").Value = 1875
$ws.Range("J68").Value = 1899.5
$ws.Range("K68").Value = 5625
$ws.Range("L68").Value = 5698.5
$ws.Range("M68").Value = -4814
$ws.Range("N68").Value = -7320.5
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H71").Value = 1883.1666
$ws.Range("I71").Value = 1875
$ws.Range("J71").Value = 1899.5
$ws.Range("K71").Value = 16875
$ws.Range("L71").Value = 17095.5
$ws.Range("M71").Value = -12819
$ws.Range("N71").Value = -25207.5
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H86").Value = 5500001.5
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 5500001.5
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H116").Value = 3016
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 3016
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 9048
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -15932
$ws.Range("H128").Value = 276308.66
$ws.Range("I128").Value = 276308.66
$ws.Range("K128").Value = 828925.98
$ws.Range("M128").Value = -823945.98

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 171.41667
$ws.Range("I2").Value = 128
$ws.Range("J2").Value = 649
$ws.Range("K2").Value = 128
$ws.Range("L2").Value = 649
$ws.Range("M2").Value = -15
$ws.Range("N2").Value = -875
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 326185.38
$ws.Range("I132").Value = 387528.78
$ws.Range("J132").Value = 7199.8
$ws.Range("K132").Value = 1162586.34
$ws.Range("L132").Value = 21599.4
$ws.Range("M132").Value = -1160056.34
$ws.Range("N132").Value = -26659.4
$ws.Range("H134").Value = 99750
$ws.Range("J134").Value = 99750
$ws.Range("L134").Value = 299250
$ws.Range("N134").Value = -304320

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1462.2273
$ws.Range("J22").Value = 2112
$ws.Range("L22").Value = 2112
$ws.Range("N22").Value = -2702
$ws.Range("H27").Value = 1462.2273
$ws.Range("J27").Value = 2112
$ws.Range("L27").Value = 2112
$ws.Range("N27").Value = -2326
$ws.Range("H136").Value = 4454.6665
$ws.Range("I136").Value = 4074.125
$ws.Range("K136").Value = 12222.375
$ws.Range("M136").Value = -9672.375
